$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rounded values for row 5 (columns B..AH)
$row5 = [ordered]@{
    "B5"  = 9.99
    "C5"  = 7.51
    "D5"  = 0.35
    "E5"  = 21.59
    "F5"  = 17.96
    "G5"  = 7.62
    "H5"  = 32.73
    "I5"  = 12.41
    "J5"  = 5.83
    "K5"  = 8.19
    "L5"  = 8.800000000000001
    "M5"  = 9.140000000000001
    "N5"  = 2.53
    "O5"  = 8
    "P5"  = 11.03
    "Q5"  = 6.28
    "R5"  = 0.6899999999999999
    "S5"  = 0.29
    "T5"  = 112.66
    "U5"  = 21.85
    "V5"  = 7.39
    "W5"  = 15.17
    "X5"  = 8.08
    "Y5"  = 1.04
    "Z5"  = 15.78
    "AA5" = 6.43
    "AB5" = 5.61
    "AC5" = 6.67
    "AD5" = 9.24
    "AE5" = 0.17
    "AF5" = 30.19
    "AG5" = 4.3
    "AH5" = 8.960000000000001
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove the last data row (row 6) entirely, shrinking the used range to A1:AH5
$ws.Rows.Item(6).Delete()
